# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row (data rows 2..38 in this sheet)
$lastRow = $ws.UsedRange.Rows.Count

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style from an existing header cell (A1) to the new headers
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wins = 97
$losses = 65
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}
